$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column letter -> new value
$changes = @{
    6  = @{ 'S' = 1.8;  'T' = 2 }
    7  = @{
        'G' = 1.62; 'H' = 3.45; 'I' = 5.7;  'J' = 2.15; 'L' = 5.7;
        'N' = 6.5;  'O' = 1.37; 'P' = 2.85; 'S' = 2.1;  'W' = 3.55;
        'X' = 1.25; 'Z' = 2.67; 'AD' = 6.9; 'AF' = 11.75; 'AG' = 14;
        'AI' = 6.5; 'AJ' = 6.8; 'AK' = 19;  'AN' = 12.5; 'AO' = 32;
        'AP' = 18;  'AQ' = 120; 'AR' = 70;  'AS' = 75
    }
    9  = @{
        'G' = 1.42; 'H' = 3.9;  'J' = 1.91; 'K' = 2.18; 'L' = 7.8;
        'N' = 6.4;  'O' = 1.38; 'P' = 2.8;  'S' = 2.12; 'T' = 1.65;
        'W' = 3.65; 'X' = 1.24; 'Z' = 2.62; 'AA' = 2.4; 'AB' = 1.5;
        'AD' = 5.4; 'AF' = 8.5; 'AI' = 6.4; 'AJ' = 8;   'AK' = 26;
        'AN' = 17;  'AP' = 28;  'AQ' = 300; 'AR' = 150
    }
    11 = @{ 'S' = 1.57; 'T' = 2.35; 'U' = 1.95; 'V' = 1.85; 'W' = 2.38; 'X' = 1.53 }
    13 = @{
        'H' = 4;   'I' = 1.57; 'J' = 5.5; 'L' = 2.1; 'O' = 1.29;
        'P' = 3.5; 'AA' = 1.91; 'AB' = 1.8; 'AJ' = 8; 'AK' = 19;
        'AQ' = 11; 'AS' = 29
    }
    15 = @{ 'I' = 2 }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        $ws.Range($addr).Value = $cols[$col]
    }
}
